# Insert a new weekly price-report row for "Femacal de La Calera" (Acelga)
# at sheet row 377, pushing the existing row 377..501 down to 378..502.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 377 (shifts 377:501 -> 378:502,
# and extends the used range to A1:R502, matching the target dimension).
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new record's data.
$ws.Cells.Item(377, 1).Value = 3
$ws.Cells.Item(377, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(377, 3).Value = "Coquimbo"
$ws.Cells.Item(377, 4).Value = "2023-02-28"
$ws.Cells.Item(377, 5).Value = 5
$ws.Cells.Item(377, 6).Value = 100112009
$ws.Cells.Item(377, 7).Value = "Acelga"
$ws.Cells.Item(377, 8).Value = "Sin especificar"
$ws.Cells.Item(377, 9).Value = "Primera"
$ws.Cells.Item(377, 10).Value = 248
$ws.Cells.Item(377, 11).Value = 3500
$ws.Cells.Item(377, 12).Value = 3800
$ws.Cells.Item(377, 13).Value = 3633
$ws.Cells.Item(377, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(377, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(377, 16).Value = 606
$ws.Cells.Item(377, 17).Value = 6
$ws.Cells.Item(377, 18).Value = "Hortaliza"
